$d = $word.ActiveDocument

# 1. Sanctioned amount in number -> prefix with rupee symbol
$d.Content.Find.Execute("{{kfsData.terms.sanctionLimit.inNumber}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{kfsData.terms.sanctionLimit.inNumber}}", 2)

# 2. Due day placeholder moved from installmentDetails to termLoanDetails (this specific
#    occurrence is the only one still referencing installmentDetails.dueDay)
$d.Content.Find.Execute("{{kfsData.installmentDetails.dueDay}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{kfsData.termLoanDetails.dueDay}}", 2)

# 3. Dishonour / cheque bounce charges -> simple variable placeholder
$d.Content.Find.Execute("{{kfsData.applicableFees.defaultingFee.baseFee.value}} + 18% GST", $true, $false, $false, $false, $false, $true, 1, $false, "{{dishonourCharges}}", 2)

# 4. Mandate swap charges -> simple variable placeholder
$d.Content.Find.Execute("{{kfsData.applicableFees.mandateSwapCharges.baseFee.value}} + 18% GST", $true, $false, $false, $false, $false, $true, 1, $false, "{{bankMandateSwapCharges}}", 2)

# 5. Renewal fee -> simple variable placeholder
$d.Content.Find.Execute("{{kfsData.applicableFees.renewalFee.baseFee.value}} + 18% GST", $true, $false, $false, $false, $false, $true, 1, $false, "{{renewalFee}}", 2)

# 6. Additional security pledging fee -> simple variable placeholder, and add a trailing
#    empty run (matching the formatting pattern used elsewhere in the document).
$rng = $d.Content
$rng.Find.Execute("{{kfsData.applicableFees.additionalSecurityPledgingFee.baseFee.value}} + 18% GST", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">{{additionalSecurityPledging}}</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# 7. Credit limit in number -> prefix with rupee symbol
$d.Content.Find.Execute("{{kfsData.terms.creditLimit.inNumber}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{kfsData.terms.creditLimit.inNumber}}", 2)

# 8. Fix missing closing braces typo on the days placeholder
$d.Content.Find.Execute("{{kfsData.terms.months}} months {{kfsData.terms.days} days", $true, $false, $false, $false, $false, $true, 1, $false, "{{kfsData.terms.months}} months {{kfsData.terms.days}} days", 2)

# 9. Total interest paid (no currency symbol occurrence) -> simple variable placeholder
$d.Content.Find.Execute("{{kfsData.installmentDetails.totalInterestPaid}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{totalInterestPaid}}", 2)

# 10. Processing fee total (appears twice) -> rupee-prefixed simple variable placeholder
$d.Content.Find.Execute("{{kfsData.applicableFees.processingFees.totalFee}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{processingFeeTotalValue}}", 2)

# 11. Net disbursement amount -> prefix with rupee symbol
$d.Content.Find.Execute("{{kfsData.installmentDetails.netDisbursementAmount}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{kfsData.installmentDetails.netDisbursementAmount}}", 2)

# 12. Total installment amount (no currency symbol occurrence) -> simple variable placeholder
$d.Content.Find.Execute("{{kfsData.installmentDetails.totalInstallmentAmount}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{totalInstallmentAmount}}", 2)

# 13. Table grid column widths (repayment schedule table)
$t = $d.Tables.Item(4)
$t.Columns.Item(5).Width = 74.25
$t.Columns.Item(6).Width = 75.0

# 14. Total installment amount (already rupee-prefixed occurrence) -> simple variable placeholder
$d.Content.Find.Execute("₹{{kfsData.installmentDetails.totalInstallmentAmount}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{totalInstallmentAmount}}", 2)

# 15. Total principal paid (rupee-prefixed occurrence) -> simple variable placeholder
$d.Content.Find.Execute("₹{{kfsData.installmentDetails.totalPrincipalPaid}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{totalPrinciplePaid}}", 2)

# 16. Total interest paid (rupee-prefixed occurrence) -> simple variable placeholder
$d.Content.Find.Execute("₹{{kfsData.installmentDetails.totalInterestPaid}}", $true, $false, $false, $false, $false, $true, 1, $false, "₹{{totalInterestPaid}}", 2)

Write-Output "done"
